$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bug fixes on existing row 5 (Laci Ferenczi): nickname + career 180s count
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Parittyás"
$ws.Range("H5").Value = 2

# ---------------------------------------------------------------------------
# 2) Propagate row 14's current ("last row") formatting down to the new row 15
#    before we touch row 14's own formatting.
# ---------------------------------------------------------------------------
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null

$ws.Range("H14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null

$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null

# E15/F15/G15 get the plain body font plus a fresh bottom border.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E15:G15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15:G15").Borders.Item(9).Color = 0
$ws.Range("E15:G15").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 3) Fill in the new row 15 data (2025 round 5 addition: Feri Gyulai-Nagy)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Feri Gyulai-Nagy"
$ws.Range("B15").Value = "Feri"
$ws.Range("C15").Value = "Gwen Stefani, Akon - The Sweet Escape"
$ws.Range("D15").Value = 39597
$ws.Range("D15").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("E15").Value = "HUN"
$ws.Range("F15").Value = "Szeged"
$ws.Range("G15").Value = "N/A"
$ws.Range("H15").Value = 0

# ---------------------------------------------------------------------------
# 4) Row 14 (Balázs Pápai) is no longer the last row - drop its bottom border
#    and re-apply the plain "middle of table" formatting.
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null

$ws.Range("C13").Copy() | Out-Null
$ws.Range("C14:G14").PasteSpecial(-4122) | Out-Null

$ws.Range("B14").Borders.Item(7).Color = 0
$ws.Range("B14").Borders.Item(7).LineStyle = 1

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = 36552
$ws.Range("D14").NumberFormat = "yyyy\-mm\-dd"

# ---------------------------------------------------------------------------
# 5) Move the active selection the way the author left it after the edit.
# ---------------------------------------------------------------------------
$ws.Range("A17").Select() | Out-Null
